# Update "想去人数" (column F) figures on the "展览" and "全部类型" sheets
# to reflect the latest scrape output (gh-pages output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# Map of row number -> new value for column F, shared by both affected sheets.
$updates = @{
    2  = 1429
    3  = 7771
    9  = 6111
    10 = 164
    11 = 17
    13 = 1847
    14 = 1379
    16 = 866
    17 = 180
    18 = 5550
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
